# Auto-generated script applying the Kujata_Profits market-data refresh.
# For each affected sheet/cell, write the updated currentAveragePrice / LevePrice / LeveProfit
# figures that the scheduled market-data runner produced.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 388  # H5
$ws.Cells.Item(5, 9).Value = 70  # I5
$ws.Cells.Item(5, 11).Value = 70  # K5
$ws.Cells.Item(5, 13).Value = 45  # M5
$ws.Cells.Item(21, 8).Value = 900  # H21
$ws.Cells.Item(21, 9).Value = 900  # I21
$ws.Cells.Item(21, 11).Value = 900  # K21
$ws.Cells.Item(21, 13).Value = -432  # M21
$ws.Cells.Item(23, 8).Value = 900  # H23
$ws.Cells.Item(23, 9).Value = 900  # I23
$ws.Cells.Item(23, 11).Value = 900  # K23
$ws.Cells.Item(23, 13).Value = -666  # M23
$ws.Cells.Item(40, 8).Value = 2025.2  # H40
$ws.Cells.Item(40, 10).Value = 1641.8889  # J40
$ws.Cells.Item(40, 12).Value = 1641.8889  # L40
$ws.Cells.Item(40, 14).Value = -1991.8889  # N40
$ws.Cells.Item(62, 8).Value = 3979.7  # H62
$ws.Cells.Item(62, 9).Value = 3966.5  # I62
$ws.Cells.Item(62, 10).Value = 3999.5  # J62
$ws.Cells.Item(62, 11).Value = 3966.5  # K62
$ws.Cells.Item(62, 12).Value = 3999.5  # L62
$ws.Cells.Item(62, 13).Value = -3342.5  # M62
$ws.Cells.Item(62, 14).Value = -5247.5  # N62
$ws.Cells.Item(64, 8).Value = 3529.7188  # H64
$ws.Cells.Item(64, 9).Value = 3522  # I64
$ws.Cells.Item(64, 10).Value = 3535  # J64
$ws.Cells.Item(64, 11).Value = 3522  # K64
$ws.Cells.Item(64, 12).Value = 3535  # L64
$ws.Cells.Item(64, 13).Value = -3274  # M64
$ws.Cells.Item(64, 14).Value = -4031  # N64
$ws.Cells.Item(65, 8).Value = 3979.7  # H65
$ws.Cells.Item(65, 9).Value = 3966.5  # I65
$ws.Cells.Item(65, 10).Value = 3999.5  # J65
$ws.Cells.Item(65, 11).Value = 19832.5  # K65
$ws.Cells.Item(65, 12).Value = 19997.5  # L65
$ws.Cells.Item(65, 13).Value = -16712.5  # M65
$ws.Cells.Item(65, 14).Value = -26237.5  # N65
$ws.Cells.Item(67, 8).Value = 3529.7188  # H67
$ws.Cells.Item(67, 9).Value = 3522  # I67
$ws.Cells.Item(67, 10).Value = 3535  # J67
$ws.Cells.Item(67, 11).Value = 3522  # K67
$ws.Cells.Item(67, 12).Value = 3535  # L67
$ws.Cells.Item(67, 13).Value = -2664  # M67
$ws.Cells.Item(67, 14).Value = -5251  # N67
$ws.Cells.Item(106, 8).Value = 3157.3  # H106
$ws.Cells.Item(106, 9).Value = 3157.3  # I106
$ws.Cells.Item(106, 10).Value = 0  # J106
$ws.Cells.Item(106, 11).Value = 3157.3  # K106
$ws.Cells.Item(106, 12).Value = 0  # L106
$ws.Cells.Item(106, 13).Value = -2526.3  # M106
$ws.Cells.Item(106, 14).ClearContents()  # N106
$ws.Cells.Item(107, 8).Value = 3849.889  # H107
$ws.Cells.Item(107, 9).Value = 3449.8333  # I107
$ws.Cells.Item(107, 10).Value = 4650  # J107
$ws.Cells.Item(107, 11).Value = 3449.8333  # K107
$ws.Cells.Item(107, 12).Value = 4650  # L107
$ws.Cells.Item(107, 13).Value = -1529.8333  # M107
$ws.Cells.Item(107, 14).Value = -8490  # N107
$ws.Cells.Item(111, 8).Value = 5682.25  # H111
$ws.Cells.Item(111, 9).Value = 5682.25  # I111
$ws.Cells.Item(111, 11).Value = 17046.75  # K111
$ws.Cells.Item(111, 13).Value = -13979.75  # M111
$ws.Cells.Item(112, 8).Value = 2058.279  # H112
$ws.Cells.Item(112, 10).Value = 2107.4634  # J112
$ws.Cells.Item(112, 12).Value = 6322.3902  # L112
$ws.Cells.Item(112, 14).Value = -8538.3902  # N112
$ws.Cells.Item(113, 8).Value = 3548.6  # H113
$ws.Cells.Item(113, 9).Value = 2445  # I113
$ws.Cells.Item(113, 10).Value = 3824.5  # J113
$ws.Cells.Item(113, 11).Value = 2445  # K113
$ws.Cells.Item(113, 12).Value = 3824.5  # L113
$ws.Cells.Item(113, 13).Value = 809  # M113
$ws.Cells.Item(113, 14).Value = -10332.5  # N113
$ws.Cells.Item(137, 8).Value = 1474.5714  # H137
$ws.Cells.Item(137, 9).Value = 1487  # I137
$ws.Cells.Item(137, 11).Value = 4461  # K137
$ws.Cells.Item(137, 13).Value = -1911  # M137
$ws.Cells.Item(138, 8).Value = 597527.7  # H138
$ws.Cells.Item(138, 9).Value = 1019.8  # I138
$ws.Cells.Item(138, 10).Value = 692211.4399999999  # J138
$ws.Cells.Item(138, 11).Value = 3059.4  # K138
$ws.Cells.Item(138, 12).Value = 2076634.32  # L138
$ws.Cells.Item(138, 13).Value = 2080.6  # M138
$ws.Cells.Item(138, 14).Value = -2086914.32  # N138

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2880.1707  # H32
$ws.Cells.Item(32, 9).Value = 3012.4866  # I32
$ws.Cells.Item(32, 10).Value = 1656.25  # J32
$ws.Cells.Item(32, 11).Value = 3012.4866  # K32
$ws.Cells.Item(32, 12).Value = 1656.25  # L32
$ws.Cells.Item(32, 13).Value = -2725.4866  # M32
$ws.Cells.Item(32, 14).Value = -2230.25  # N32
$ws.Cells.Item(61, 8).Value = 1415.7778  # H61
$ws.Cells.Item(61, 9).Value = 1217.75  # I61
$ws.Cells.Item(61, 10).Value = 3000  # J61
$ws.Cells.Item(61, 11).Value = 1217.75  # K61
$ws.Cells.Item(61, 12).Value = 3000  # L61
$ws.Cells.Item(61, 13).Value = -1005.75  # M61
$ws.Cells.Item(61, 14).Value = -3424  # N61
$ws.Cells.Item(104, 8).Value = 0  # H104
$ws.Cells.Item(104, 10).Value = 0  # J104
$ws.Cells.Item(104, 12).Value = 0  # L104
$ws.Cells.Item(104, 14).ClearContents()  # N104
$ws.Cells.Item(122, 8).Value = 1766.6666  # H122
$ws.Cells.Item(122, 9).Value = 1800  # I122
$ws.Cells.Item(122, 10).Value = 1700  # J122
$ws.Cells.Item(122, 11).Value = 5400  # K122
$ws.Cells.Item(122, 12).Value = 5100  # L122
$ws.Cells.Item(122, 13).Value = -2950  # M122
$ws.Cells.Item(122, 14).Value = -10000  # N122
$ws.Cells.Item(132, 8).Value = 2364.647  # H132
$ws.Cells.Item(132, 9).Value = 2050.0356  # I132
$ws.Cells.Item(132, 11).Value = 6150.1068  # K132
$ws.Cells.Item(132, 13).Value = -3620.1068  # M132
$ws.Cells.Item(136, 8).Value = 1415.7778  # H136
$ws.Cells.Item(136, 9).Value = 1217.75  # I136
$ws.Cells.Item(136, 10).Value = 3000  # J136
$ws.Cells.Item(136, 11).Value = 3653.25  # K136
$ws.Cells.Item(136, 12).Value = 9000  # L136
$ws.Cells.Item(136, 13).Value = -1103.25  # M136
$ws.Cells.Item(136, 14).Value = -14100  # N136

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 90910160  # H16
$ws.Cells.Item(16, 9).Value = 125001020  # I16
$ws.Cells.Item(16, 10).Value = 1166.6666  # J16
$ws.Cells.Item(16, 11).Value = 125001020  # K16
$ws.Cells.Item(16, 12).Value = 1166.6666  # L16
$ws.Cells.Item(16, 13).Value = -125000733  # M16
$ws.Cells.Item(16, 14).Value = -1740.6666  # N16
$ws.Cells.Item(31, 8).Value = 1013.125  # H31
$ws.Cells.Item(31, 9).Value = 920.6667  # I31
$ws.Cells.Item(31, 11).Value = 920.6667  # K31
$ws.Cells.Item(31, 13).Value = -625.6667  # M31
$ws.Cells.Item(34, 8).Value = 1013.125  # H34
$ws.Cells.Item(34, 9).Value = 920.6667  # I34
$ws.Cells.Item(34, 11).Value = 920.6667  # K34
$ws.Cells.Item(34, 13).Value = -718.6667  # M34
$ws.Cells.Item(113, 8).Value = 90910160  # H113
$ws.Cells.Item(113, 9).Value = 125001020  # I113
$ws.Cells.Item(113, 10).Value = 1166.6666  # J113
$ws.Cells.Item(113, 11).Value = 125001020  # K113
$ws.Cells.Item(113, 12).Value = 1166.6666  # L113
$ws.Cells.Item(113, 13).Value = -124998850  # M113
$ws.Cells.Item(113, 14).Value = -5506.6666  # N113
$ws.Cells.Item(132, 8).Value = 9334.066000000001  # H132
$ws.Cells.Item(132, 9).Value = 10346.917  # I132
$ws.Cells.Item(132, 10).Value = 5282.6665  # J132
$ws.Cells.Item(132, 11).Value = 31040.751  # K132
$ws.Cells.Item(132, 12).Value = 15847.9995  # L132
$ws.Cells.Item(132, 13).Value = -28510.751  # M132
$ws.Cells.Item(132, 14).Value = -20907.9995  # N132

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 301127.66  # H4
$ws.Cells.Item(4, 9).Value = 299729.66  # I4
$ws.Cells.Item(4, 10).Value = 301508.9  # J4
$ws.Cells.Item(4, 11).Value = 899188.98  # K4
$ws.Cells.Item(4, 12).Value = 904526.7000000001  # L4
$ws.Cells.Item(4, 13).Value = -899076.98  # M4
$ws.Cells.Item(4, 14).Value = -904750.7000000001  # N4
$ws.Cells.Item(22, 8).Value = 2900  # H22
$ws.Cells.Item(22, 10).Value = 2600  # J22
$ws.Cells.Item(22, 12).Value = 7800  # L22
$ws.Cells.Item(22, 14).Value = -8138  # N22
$ws.Cells.Item(27, 8).Value = 2900  # H27
$ws.Cells.Item(27, 10).Value = 2600  # J27
$ws.Cells.Item(27, 12).Value = 7800  # L27
$ws.Cells.Item(27, 14).Value = -8004  # N27
$ws.Cells.Item(41, 8).Value = 398  # H41
$ws.Cells.Item(41, 9).Value = 398  # I41
$ws.Cells.Item(41, 10).Value = 0  # J41
$ws.Cells.Item(41, 11).Value = 1194  # K41
$ws.Cells.Item(41, 12).Value = 0  # L41
$ws.Cells.Item(41, 13).Value = -856  # M41
$ws.Cells.Item(41, 14).ClearContents()  # N41
$ws.Cells.Item(44, 8).Value = 1633.3334  # H44
$ws.Cells.Item(44, 9).Value = 700  # I44
$ws.Cells.Item(44, 11).Value = 2100  # K44
$ws.Cells.Item(44, 13).Value = -1702  # M44
$ws.Cells.Item(139, 8).Value = 3434.5454  # H139
$ws.Cells.Item(139, 9).Value = 3438.5557  # I139
$ws.Cells.Item(139, 10).Value = 3416.5  # J139
$ws.Cells.Item(139, 11).Value = 10315.6671  # K139
$ws.Cells.Item(139, 12).Value = 10249.5  # L139
$ws.Cells.Item(139, 13).Value = -5175.667099999999  # M139
$ws.Cells.Item(139, 14).Value = -20529.5  # N139

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2387.4482  # H132
$ws.Cells.Item(132, 9).Value = 2049.52  # I132
$ws.Cells.Item(132, 11).Value = 6148.559999999999  # K132
$ws.Cells.Item(132, 13).Value = -3618.559999999999  # M132

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 1050004  # H25
$ws.Cells.Item(25, 10).Value = 100008  # J25
$ws.Cells.Item(25, 12).Value = 100008  # L25
$ws.Cells.Item(25, 14).Value = -100468  # N25
$ws.Cells.Item(40, 8).Value = 3068  # H40
$ws.Cells.Item(40, 9).Value = 2597.1428  # I40
$ws.Cells.Item(40, 11).Value = 2597.1428  # K40
$ws.Cells.Item(40, 13).Value = -2461.1428  # M40
$ws.Cells.Item(46, 8).Value = 7916.75  # H46
$ws.Cells.Item(46, 9).Value = 1500.3334  # I46
$ws.Cells.Item(46, 11).Value = 1500.3334  # K46
$ws.Cells.Item(46, 13).Value = -1312.3334  # M46
$ws.Cells.Item(61, 8).Value = 2261.6667  # H61
$ws.Cells.Item(61, 9).Value = 1980  # I61
$ws.Cells.Item(61, 10).Value = 2402.5  # J61
$ws.Cells.Item(61, 11).Value = 1980  # K61
$ws.Cells.Item(61, 12).Value = 2402.5  # L61
$ws.Cells.Item(61, 13).Value = -1778  # M61
$ws.Cells.Item(61, 14).Value = -2806.5  # N61
$ws.Cells.Item(82, 8).Value = 1901.6666  # H82
$ws.Cells.Item(82, 9).Value = 1302  # I82
$ws.Cells.Item(82, 10).Value = 2201.5  # J82
$ws.Cells.Item(82, 11).Value = 1302  # K82
$ws.Cells.Item(82, 12).Value = 2201.5  # L82
$ws.Cells.Item(82, 13).Value = -941  # M82
$ws.Cells.Item(82, 14).Value = -2923.5  # N82
$ws.Cells.Item(85, 8).Value = 1901.6666  # H85
$ws.Cells.Item(85, 9).Value = 1302  # I85
$ws.Cells.Item(85, 10).Value = 2201.5  # J85
$ws.Cells.Item(85, 11).Value = 1302  # K85
$ws.Cells.Item(85, 12).Value = 2201.5  # L85
$ws.Cells.Item(85, 13).Value = -54  # M85
$ws.Cells.Item(85, 14).Value = -4697.5  # N85
$ws.Cells.Item(100, 8).Value = 2166.6667  # H100
$ws.Cells.Item(113, 8).Value = 2261.6667  # H113
$ws.Cells.Item(113, 9).Value = 1980  # I113
$ws.Cells.Item(113, 10).Value = 2402.5  # J113
$ws.Cells.Item(113, 11).Value = 1980  # K113
$ws.Cells.Item(113, 12).Value = 2402.5  # L113
$ws.Cells.Item(113, 13).Value = 190  # M113
$ws.Cells.Item(113, 14).Value = -6742.5  # N113
$ws.Cells.Item(132, 8).Value = 34135.484  # H132
$ws.Cells.Item(132, 9).Value = 1628.381  # I132
$ws.Cells.Item(132, 10).Value = 102400.4  # J132
$ws.Cells.Item(132, 11).Value = 4885.143  # K132
$ws.Cells.Item(132, 12).Value = 307201.2  # L132
$ws.Cells.Item(132, 13).Value = -2355.143  # M132
$ws.Cells.Item(132, 14).Value = -312261.2  # N132

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 393.125  # H100
$ws.Cells.Item(100, 9).Value = 393.125  # I100
$ws.Cells.Item(100, 11).Value = 786.25  # K100
$ws.Cells.Item(100, 13).Value = -245.25  # M100
$ws.Cells.Item(108, 8).Value = 3026  # H108
$ws.Cells.Item(108, 10).Value = 3026  # J108
$ws.Cells.Item(108, 12).Value = 3026  # L108
$ws.Cells.Item(108, 14).Value = -10706  # N108
$ws.Cells.Item(113, 8).Value = 510.94116  # H113
$ws.Cells.Item(113, 9).Value = 390.36365  # I113
$ws.Cells.Item(113, 11).Value = 1171.09095  # K113
$ws.Cells.Item(113, 13).Value = 998.90905  # M113
$ws.Cells.Item(133, 8).Value = 34931.668  # H133
$ws.Cells.Item(133, 10).Value = 34931.668  # J133
$ws.Cells.Item(133, 11).Value = 34931.668  # K133
$ws.Cells.Item(133, 12).Value = 34931.668  # L133
$ws.Cells.Item(133, 14).Value = -45051.668  # N133
